# Add LAFC to team name links
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A43").Value = "LAFC"
$ws.Range("B43").Value = "LAFC"

$ws.Range("A44").Value = "Los Angeles FC"
$ws.Range("B44").Value = "LAFC"

$ws.Range("C37").Select()
